$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column (C) for rows 2-15 from
# serial date 45224 (2023-10-25) to 45233 (2023-11-03).
for ($row = 2; $row -le 15; $row++) {
    $ws.Cells.Item($row, 3).Value = 45233
}
